$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Update status text "In Translation" -> "Ready for handoff" on all sheets
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("C2").Value = "Ready for handoff"

# Update the handoff generation timestamps
$wsOverview.Range("G2").Value = "2016-08-28 10:47:37"
$wsDeDe.Range("H2").Value = "2016-08-28 10:47:37"
$wsZhCn.Range("H2").Value = "2016-08-28 10:47:33"

# Widen the affected status columns to account for the longer text
$wsOverview.Columns.Item(5).ColumnWidth = 16.3333333333333
$wsOverview.Columns.Item(6).ColumnWidth = 16.3333333333333
$wsZhCn.Columns.Item(3).ColumnWidth = 16.3333333333333
$wsDeDe.Columns.Item(3).ColumnWidth = 16.3333333333333
